$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two header cells (text content stays associated with same column,
# but label text is updated to reflect km/hrs instead of miles/hours).
$ws.Range("D1").Value = "Time Duration (hrs)"
$ws.Range("C1").Value = "Distance Travelled (kms)"

# Add a sum total row under the Distance column.
$ws.Range("C15").Formula = "=SUM(C3:C14)"
$ws.Range("C15").Font.Bold = $true

# Adjust column widths to match new content.
$ws.Columns.Item(3).ColumnWidth = 20.833333333333332
$ws.Columns.Item(4).ColumnWidth = 16.5
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668

# Update the view: select F16 (below the new total row).
$ws.Range("F16").Select()
